$d = $word.ActiveDocument

# 1. Fix the "input da ta into" -> "input data into" typo (a stray space was
#    left splitting "data" across runs; collapse it back into one word).
$d.Content.Find.Execute("input da ta into the console", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "input data into the console", 2)

# 2. Move the "_GoBack" bookmark from the end of the "Description of
#    Understanding" paragraph to the end of the "Teaching Video: None"
#    paragraph (it simply tracks the last edit location in the document).
if ($d.Bookmarks.Exists("_GoBack")) {
    $old = $d.Bookmarks("_GoBack")
    $oldStart = $old.Start
    $old.Delete()
}

$tv = $d.Content.Find
$found = $d.Content.Find.Execute("Teaching Video:" + [char]13 + [char]10, $false, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)

# Locate the paragraph that starts with "Teaching Video:" and place the new
# bookmark right after its text (" None"), before the paragraph mark.
foreach ($p in $d.Paragraphs) {
    $ptext = $p.Range.Text
    if ($ptext.StartsWith("Teaching Video:") -and $ptext.Contains("None")) {
        $pEnd = $p.Range.End - 1
        $target = $d.Range($pEnd, $pEnd)
        $d.Bookmarks.Add("_GoBack", $target)
        break
    }
}
